$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the cryptos price/volume table (GitHub Actions job).
# Price/Volume columns are stored as plain text in the source data, so for
# any new value that looks numeric we briefly force a text number format
# before writing it -- otherwise Excel would auto-convert it to a number
# (losing formatting like trailing zeros / thousands separators) -- then
# restore the cells original style.

$ws.Range("D2").Value = '61.444.14'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '2.992.83'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '599.05'
$c.Style = $s
$ws.Range("E5").Value = '  +3.05%  '
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '143.97'
$c.Style = $s
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").Value = '2.990.84'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  -0.94%  '
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.07'
$c.Style = $s
$ws.Range("E11").Value = '  +7.49%  '
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("E13").Value = '  +0.52%  '
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '34.39'
$c.Style = $s
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D16").Value = '3.488.42'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").Value = '61.437.65'
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").Value = '2.990.55'
$ws.Range("E19").Value = '  -0.36%  '
$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '450.81'
$c.Style = $s
$ws.Range("E20").Value = '  -0.76%  '
$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '14.01'
$c.Style = $s
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("E22").Value = '  +1.56%  '
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.33'
$c.Style = $s
$ws.Range("E23").Value = '  +0.58%  '
$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '81.53'
$c.Style = $s
$ws.Range("E24").Value = '  +2.07%  '
$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '10.76'
$c.Style = $s
$ws.Range("E25").Value = '  +6.42%  '
$ws.Range("E26").Value = '  -3.21%  '
$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.94'
$c.Style = $s
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("E28").Value = '  +0.22%  '
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.Style = $s
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("E32").Value = '  -1.49%  '
$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.26'
$c.Style = $s
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("E34").Value = '  +3.39%  '
$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  +5.06%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  +1.04%  '
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.21'
$c.Style = $s
$ws.Range("E38").Value = '  +3.14%  '
$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '50.40'
$c.Style = $s
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("E41").Value = '  +11.18%  '
$ws.Range("E42").Value = '  -0.86%  '
$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '397.54'
$c.Style = $s
$ws.Range("E43").Value = '  -3.37%  '
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '39.96'
$c.Style = $s
$ws.Range("E44").Value = '  +5.31%  '
$ws.Range("E45").Value = '  +0.50%  '
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.271'
$c.Style = $s
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").Value = '2.688.95'
$ws.Range("E47").Value = '  -2.61%  '
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '131.04'
$c.Style = $s
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  -0.50%  '
$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.15'
$c.Style = $s
$ws.Range("E51").Value = '  +1.67%  '
